# Auto update Excel log 2026-02-04 14:21:25
# Appends new sensor log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($range, [string]$text) {
    # Force the cell to be treated as plain text so date/number-like
    # strings (e.g. "2026-02-04") are not auto-converted by Excel's
    # input parser, then drop the temporary Text format so the cell
    # is left with the default (General) style, matching the rest of
    # the log.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Add-LogRows($sheetName, $rows) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $rows) {
        $r = $row[0]
        Set-TextCell $ws.Cells.Item($r, 1) $row[1]
        Set-TextCell $ws.Cells.Item($r, 2) $row[2]
        Set-TextCell $ws.Cells.Item($r, 3) $row[3]
        Set-TextCell $ws.Cells.Item($r, 4) $row[4]
        Set-TextCell $ws.Cells.Item($r, 5) $row[5]
        Set-TextCell $ws.Cells.Item($r, 6) $row[6]
    }
}

# PIR sheet: rows 231-243 (Date, Timestamp, Hour, Location, Value, Status)
$pirRows = @(
    @(231, "2026-02-04", "14:20:21", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(232, "2026-02-04", "14:20:23", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(233, "2026-02-04", "14:20:27", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(234, "2026-02-04", "14:20:32", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(235, "2026-02-04", "14:20:37", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(236, "2026-02-04", "14:20:42", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(237, "2026-02-04", "14:20:47", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(238, "2026-02-04", "14:20:52", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(239, "2026-02-04", "14:20:57", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(240, "2026-02-04", "14:21:02", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(241, "2026-02-04", "14:21:07", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(242, "2026-02-04", "14:21:12", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(243, "2026-02-04", "14:21:17", "14:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows "PIR" $pirRows

# Humidity sheet: rows 195-202
$humidityRows = @(
    @(195, "2026-02-04", "14:20:22", "14:00", "Bathroom", "78.9%", "Active"),
    @(196, "2026-02-04", "14:20:23", "14:00", "Bathroom", "78.0%", "Active"),
    @(197, "2026-02-04", "14:20:28", "14:00", "Bathroom", "79.0%", "Active"),
    @(198, "2026-02-04", "14:20:33", "14:00", "Bathroom", "78.0%", "Active"),
    @(199, "2026-02-04", "14:20:38", "14:00", "Bathroom", "78.7%", "Active"),
    @(200, "2026-02-04", "14:20:43", "14:00", "Bathroom", "77.6%", "Active"),
    @(201, "2026-02-04", "14:20:59", "14:00", "Bathroom", "77.8%", "Active"),
    @(202, "2026-02-04", "14:21:09", "14:00", "Bathroom", "77.5%", "Active")
)
Add-LogRows "Humidity" $humidityRows

# Temperature sheet: rows 195-202
$temperatureRows = @(
    @(195, "2026-02-04", "14:20:22", "14:00", "Bathroom", "24.5C", "Active"),
    @(196, "2026-02-04", "14:20:24", "14:00", "Bathroom", "24.5C", "Active"),
    @(197, "2026-02-04", "14:20:29", "14:00", "Bathroom", "24.5C", "Active"),
    @(198, "2026-02-04", "14:20:34", "14:00", "Bathroom", "24.6C", "Active"),
    @(199, "2026-02-04", "14:20:39", "14:00", "Bathroom", "24.5C", "Active"),
    @(200, "2026-02-04", "14:20:44", "14:00", "Bathroom", "24.6C", "Active"),
    @(201, "2026-02-04", "14:20:59", "14:00", "Bathroom", "24.6C", "Active"),
    @(202, "2026-02-04", "14:21:09", "14:00", "Bathroom", "24.6C", "Active")
)
Add-LogRows "Temperature" $temperatureRows
